$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New random-run data for column A (B stays constant at 544732)
$newValues = @{
    2  = 17403600
    3  = 12462300
    4  = 1357300
    5  = 2256200
    6  = 1238000
    7  = 1308200
    8  = 1259400
    9  = 1097800
    10 = 1520300
    11 = 5142600
    12 = 10183100
    13 = 3272200
    14 = 7062100
    15 = 8844100
    16 = 4281200
    17 = 3002000
    18 = 1329300
    19 = 2981000
    20 = 10080000
    21 = 6896100
    22 = 1118000
    23 = 1167200
    24 = 1215800
    25 = 1132600
    26 = 11329600
    27 = 7597300
    28 = 6707200
    29 = 1600500
    30 = 6063900
    31 = 885600
    32 = 896700
    33 = 7938700
    34 = 903400
    35 = 1135100
    36 = 858400
    37 = 901800
    38 = 889100
    39 = 879200
    40 = 871600
    41 = 889500
    42 = 908400
    43 = 875400
    44 = 874400
    45 = 1086700
    46 = 1123000
    47 = 3740200
    48 = 1452000
    49 = 1528600
    50 = 1374700
    51 = 1552100
}

foreach ($row in $newValues.Keys) {
    $ws.Range("A$row").Value = $newValues[$row]
}

# Replace the static sum/average values in E2/F2 and E4/F4 with live formulas
$ws.Range("E2").Formula = "=AVERAGE(A2:A51)"
$ws.Range("F2").Formula = "=AVERAGE(B2:B51)"
$ws.Range("E4").Formula = "=E2/10000"
$ws.Range("F4").Formula = "=F2/10000"

# Restore the active selection to F2 (as seen after regenerating results)
$ws.Range("F2").Select()

$wb.Save()
